# Apply updated KNN-imputed values to result_data_KNN worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = -8.01
    "B7"  = 5.697000000000001
    "A10" = -21.604
    "A12" = -21.606
    "B15" = 5.072000000000001
    "A18" = -22.002
    "D18" = -8.638
    "D19" = -7.912000000000001
    "B20" = 7.633
    "D27" = -8.358000000000001
    "B29" = 5.645
    "B30" = 5.694000000000001
    "B31" = 6.508
    "A37" = -19.92
    "B40" = 9.327999999999999
    "D42" = -8.258000000000001
    "D44" = -8.125999999999999
    "D47" = -7.617999999999999
    "A55" = -21.942
    "D58" = -8.370999999999999
    "A68" = -21.50700000000001
    "B68" = 5.881
    "D73" = -8.214000000000002
    "B76" = 6.341000000000001
    "A77" = -20.637
    "A78" = -20.126
    "B87" = 4.462
    "B88" = 4.853999999999999
    "D95" = -8.011999999999999
    "B96" = 6.469000000000001
    "B98" = 5.828
    "B101" = 8.193999999999999
    "D101" = -8.16
    "B102" = 7.739999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
